$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 29.66705863408808
$ws.Range("D2").Value = 0.5170586340880767
$ws.Range("E2").Value = 0.2673496310850276
$ws.Range("C3").Value = 30.15988454009383
$ws.Range("D3").Value = 0.8098845400938259
$ws.Range("E3").Value = 0.655912968282988
$ws.Range("C4").Value = 29.84494678620212
$ws.Range("D4").Value = 0.4749467862021177
$ws.Range("E4").Value = 0.2255744497237201
$ws.Range("C5").Value = 29.24820229833501
$ws.Range("D5").Value = -0.2917977016649935
$ws.Range("E5").Value = 0.08514589869697256
$ws.Range("C6").Value = 29.77028807683576
$ws.Range("D6").Value = 0.2202880768357574
$ws.Range("E6").Value = 0.04852683679599654
$ws.Range("C7").Value = 29.56727933301249
$ws.Range("D7").Value = -0.1827206669875103
$ws.Range("E7").Value = 0.03338684214436065
$ws.Range("C8").Value = 30.08935803597058
$ws.Range("D8").Value = 0.2493580359705767
$ws.Range("E8").Value = 0.06217943010310343
$ws.Range("C9").Value = 29.98617848557869
$ws.Range("D9").Value = 0.1761784855786885
$ws.Range("E9").Value = 0.03103885878080017
$ws.Range("C10").Value = 29.73711101533921
$ws.Range("D10").Value = -0.1828889846607886
$ws.Range("E10").Value = 0.03344838071025417
$ws.Range("C11").Value = 29.69799928907306
$ws.Range("D11").Value = -0.2820007109269405
$ws.Range("E11").Value = 0.07952440096329984
$ws.Range("C12").Value = 29.72613621387445
$ws.Range("D12").Value = -0.3138637861255518
$ws.Range("E12").Value = 0.0985104762410661
$ws.Range("C13").Value = 29.84922737223416
$ws.Range("D13").Value = -0.3607726277658401
$ws.Range("E13").Value = 0.1301568889450694
$ws.Range("C14").Value = 29.87118448474057
$ws.Range("D14").Value = -0.3488155152594246
$ws.Range("E14").Value = 0.1216722636856978
$ws.Range("C15").Value = 29.98279341359101
$ws.Range("D15").Value = -0.3972065864089913
$ws.Range("E15").Value = 0.1577730722866835
$ws.Range("C16").Value = 30.26108524161952
$ws.Range("D16").Value = -0.1789147583804827
$ws.Range("E16").Value = 0.03201049076634649
$ws.Range("C17").Value = 30.15278547970684
$ws.Range("D17").Value = -0.3272145202931647
$ws.Range("E17").Value = 0.1070693422906859
$ws.Range("C18").Value = 30.19514282129561
$ws.Range("D18").Value = -0.4948571787043896
$ws.Range("E18").Value = 0.2448836273152682
$ws.Range("C19").Value = 30.60131952062666
$ws.Range("D19").Value = -0.1486804793733363
$ws.Range("E19").Value = 0.02210588494668508
$ws.Range("C20").Value = 30.67738320357709
$ws.Range("D20").Value = -0.2626167964229076
$ws.Range("E20").Value = 0.06896758176343089
$ws.Range("C21").Value = 30.85919113487759
$ws.Range("D21").Value = -0.09080886512240482
$ws.Range("E21").Value = 0.008246249984819112
$ws.Range("C22").Value = 31.09757034944756
$ws.Range("D22").Value = 0.07757034944756214
$ws.Range("E22").Value = 0.006017159113416904
$ws.Range("C23").Value = 31.27601048597176
$ws.Range("D23").Value = 0.156010485971759
$ws.Range("E23").Value = 0.02433927173314442
$ws.Range("C24").Value = 31.18419340973565
$ws.Range("D24").Value = -0.09580659026435612
$ws.Range("E24").Value = 0.009178902738082216
$ws.Range("C25").Value = 31.20860949197808
$ws.Range("D25").Value = -0.1713905080219149
$ws.Range("E25").Value = 0.02937470624001008
$ws.Range("C26").Value = 31.41963735202371
$ws.Range("D26").Value = -0.1603626479762887
$ws.Range("E26").Value = 0.0257161788659671
$ws.Range("C27").Value = 31.95921673663635
$ws.Range("D27").Value = 0.3092167366363476
$ws.Range("E27").Value = 0.09561499021603237
$ws.Range("C28").Value = 32.68495774151191
$ws.Range("D28").Value = 0.8049577415119096
$ws.Range("E28").Value = 0.6479569656199542
$ws.Range("C29").Value = 32.92749544001664
$ws.Range("D29").Value = 0.6474954400166339
$ws.Range("E29").Value = 0.4192503448423344
$ws.Range("C30").Value = 33.06663339335906
$ws.Range("D30").Value = 0.6166333933590522
$ws.Range("E30").Value = 0.3802367418054996
$ws.Range("C31").Value = 33.23640948043696
$ws.Range("D31").Value = 0.38640948043696
$ws.Range("E31").Value = 0.1493122865715614
$ws.Range("C32").Value = 33.29977708230282
$ws.Range("D32").Value = 0.3997770823028262
$ws.Range("E32").Value = 0.1598217155345606
$ws.Range("C33").Value = 33.33975343964335
$ws.Range("D33").Value = 0.2397534396433443
$ws.Range("E33").Value = 0.05748171182081473
$ws.Range("C34").Value = 33.71583827279082
$ws.Range("D34").Value = 0.3158382727908204
$ws.Range("E34").Value = 0.0997538145594887
$ws.Range("C35").Value = 33.74213643898327
$ws.Range("D35").Value = 0.04213643898326325
$ws.Range("E35").Value = 0.001775479490190267
$ws.Range("C36").Value = 33.83770196888512
$ws.Range("D36").Value = -0.2622980311148808
$ws.Range("E36").Value = 0.06880025712674298
$ws.Range("C37").Value = 34.26175461750775
$ws.Range("D37").Value = -0.1382453824922436
$ws.Range("E37").Value = 0.01911178578042674
$ws.Range("C38").Value = 34.53506441515873
$ws.Range("D38").Value = -0.3649355848412696
$ws.Range("E38").Value = 0.1331779810834395
$ws.Range("C39").Value = 35.41163110728882
$ws.Range("D39").Value = 0.1116311072888223
$ws.Range("E39").Value = 0.01246150411452856
$ws.Range("C40").Value = 35.68715858908152
$ws.Range("D40").Value = -0.01284141091848312
$ws.Range("E40").Value = 0.0001649018343773375
$ws.Range("C41").Value = 36.03467250596535
$ws.Range("D41").Value = -0.2653274940346506
$ws.Range("E41").Value = 0.07039867909070756
$ws.Range("C42").Value = 36.6380964444557
$ws.Range("D42").Value = -0.1619035555442991
$ws.Range("E42").Value = 0.02621276129788596
$ws.Range("C43").Value = 36.75375374323595
$ws.Range("D43").Value = -0.54624625676405
$ws.Range("E43").Value = 0.2983849730287365
$ws.Range("C44").Value = 37.77578137501172
$ws.Range("D44").Value = -0.1242186249882806
$ws.Range("E44").Value = 0.0154302667939791
$ws.Range("C45").Value = 38.49756896558319
$ws.Range("D45").Value = -0.002431034416808586
$ws.Range("E45").Value = 0.000005909928335707863
$ws.Range("C46").Value = 39.1968811618345
$ws.Range("D46").Value = 0.2968811618344986
$ws.Range("E46").Value = 0.08813842425220177
$ws.Range("C47").Value = 39.67226180578425
$ws.Range("D47").Value = 0.2722618057842467
$ws.Range("E47").Value = 0.07412649088889889
$ws.Range("C48").Value = 39.84535288553736
$ws.Range("D48").Value = -0.05464711446263948
$ws.Range("E48").Value = 0.002986307119092821
$ws.Range("C49").Value = 40.27265069933905
$ws.Range("D49").Value = 0.1726506993390515
$ws.Range("E49").Value = 0.02980826398226357
$ws.Range("C50").Value = 40.5256658712525
$ws.Range("D50").Value = -0.07433412874750189
$ws.Range("E50").Value = 0.005525562696650188
$ws.Range("C51").Value = 40.86033174156339
$ws.Range("D51").Value = -0.03966825843661326
$ws.Range("E51").Value = 0.001573570727393939
$ws.Range("C52").Value = 0.959122392995134
$ws.Range("E52").Value = 5.465621484408994
$ws.Range("E53").Value = 0.1093124296881799
